# Auto-generated Excel COM-interop script to apply crypto price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "63.967.37"
$ws.Cells.Item(2, 5).Value = "  -3.90%  "
$ws.Cells.Item(3, 4).Value = "3.612.02"
$ws.Cells.Item(3, 5).Value = "  +0.44%  "
$ws.Cells.Item(4, 4).Value = "'1.01"
$ws.Cells.Item(4, 5).Value = "  +0.38%  "
$ws.Cells.Item(5, 4).Value = "'404.02"
$ws.Cells.Item(5, 5).Value = "  -2.56%  "
$ws.Cells.Item(6, 4).Value = "'131.88"
$ws.Cells.Item(6, 5).Value = "  +1.05%  "
$ws.Cells.Item(7, 4).Value = "3.609.41"
$ws.Cells.Item(7, 5).Value = "  +0.52%  "
$ws.Cells.Item(8, 4).Value = "'0.614"
$ws.Cells.Item(8, 5).Value = "  -5.43%  "
$ws.Cells.Item(9, 5).Value = "  +0.20%  "
$ws.Cells.Item(10, 4).Value = "'0.718"
$ws.Cells.Item(10, 5).Value = "  -6.88%  "
$ws.Cells.Item(11, 4).Value = "'0.157"
$ws.Cells.Item(11, 5).Value = "  -9.57%  "
$ws.Cells.Item(12, 4).Value = "'0.0000307"
$ws.Cells.Item(12, 5).Value = "  -5.35%  "
$ws.Cells.Item(13, 4).Value = "'41.38"
$ws.Cells.Item(14, 4).Value = "'9.84"
$ws.Cells.Item(14, 5).Value = "  -0.13%  "
$ws.Cells.Item(15, 4).Value = "4.200.54"
$ws.Cells.Item(15, 5).Value = "  +0.97%  "
$ws.Cells.Item(17, 4).Value = "3.619.85"
$ws.Cells.Item(17, 5).Value = "  +1.19%  "
$ws.Cells.Item(18, 2).Value = "Chainlink"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(18, 4).Value = "'19.78"
$ws.Cells.Item(18, 5).Value = "  -2.43%  "
$ws.Cells.Item(19, 2).Value = "Uniswap"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Cells.Item(19, 4).Value = "'13.33"
$ws.Cells.Item(19, 5).Value = "  +8.87%  "
$ws.Cells.Item(20, 4).Value = "'1.07"
$ws.Cells.Item(20, 5).Value = "  -5.92%  "
$ws.Cells.Item(21, 4).Value = "64.310.48"
$ws.Cells.Item(21, 5).Value = "  -3.42%  "
$ws.Cells.Item(22, 4).Value = "'416.54"
$ws.Cells.Item(22, 5).Value = "  -6.62%  "
$ws.Cells.Item(23, 4).Value = "'14.93"
$ws.Cells.Item(23, 5).Value = "  +14.45%  "
$ws.Cells.Item(24, 4).Value = "'84.79"
$ws.Cells.Item(24, 5).Value = "  -4.52%  "
$ws.Cells.Item(25, 4).Value = "'2.96"
$ws.Cells.Item(25, 5).Value = "  -5.59%  "
$ws.Cells.Item(26, 4).Value = "'35.30"
$ws.Cells.Item(26, 5).Value = "  +0.76%  "
$ws.Cells.Item(27, 5).Value = "  -4.67%  "
$ws.Cells.Item(28, 4).Value = "'9.31"
$ws.Cells.Item(28, 5).Value = "  -6.35%  "
$ws.Cells.Item(30, 4).Value = "'12.69"
$ws.Cells.Item(30, 5).Value = "  +3.11%  "
$ws.Cells.Item(31, 4).Value = "'2.70"
$ws.Cells.Item(31, 5).Value = "  -2.26%  "
$ws.Cells.Item(32, 4).Value = "'0.115"
$ws.Cells.Item(32, 5).Value = "  -1.52%  "
$ws.Cells.Item(33, 4).Value = "'6.87"
$ws.Cells.Item(33, 5).Value = "  -6.88%  "
$ws.Cells.Item(34, 2).Value = "Kaspa"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(34, 4).Value = "'0.159"
$ws.Cells.Item(34, 5).Value = "  -0.37%  "
$ws.Cells.Item(35, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(35, 4).Value = "'40.78"
$ws.Cells.Item(35, 5).Value = "  +2.89%  "
$ws.Cells.Item(36, 2).Value = "Dai"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(36, 4).Value = "'1.00"
$ws.Cells.Item(36, 5).Value = "  +0.03%  "
$ws.Cells.Item(37, 2).Value = "OKB"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(37, 4).Value = "'55.49"
$ws.Cells.Item(37, 5).Value = "  -1.87%  "
$ws.Cells.Item(38, 5).Value = "  -5.97%  "
$ws.Cells.Item(39, 4).Value = "'2.87"
$ws.Cells.Item(39, 5).Value = "  +26.19%  "
$ws.Cells.Item(40, 4).Value = "'0.994"
$ws.Cells.Item(40, 5).Value = "  -0.50%  "
$ws.Cells.Item(41, 5).Value = "  -5.21%  "
$ws.Cells.Item(42, 4).Value = "'3.15"
$ws.Cells.Item(42, 5).Value = "  +23.15%  "
$ws.Cells.Item(43, 4).Value = "'4.38"
$ws.Cells.Item(43, 5).Value = "  +1.62%  "
$ws.Cells.Item(44, 2).Value = "Monero"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(44, 4).Value = "'144.52"
$ws.Cells.Item(44, 5).Value = "  -2.92%  "
$ws.Cells.Item(45, 2).Value = "PEPE"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Cells.Item(45, 4).Value = "0.0₃0632"
$ws.Cells.Item(45, 5).Value = "  -11.43%  "
$ws.Cells.Item(46, 2).Value = "EnergySwap"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(46, 4).Value = "'26.44"
$ws.Cells.Item(46, 5).Value = "  +24.50%  "
$ws.Cells.Item(47, 2).Value = "LidoDAOToken"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(47, 4).Value = "'3.26"
$ws.Cells.Item(47, 5).Value = "  +0.30%  "
$ws.Cells.Item(48, 4).Value = "'2.05"
$ws.Cells.Item(48, 5).Value = "  +4.12%  "
$ws.Cells.Item(49, 5).Value = "  -6.36%  "
$ws.Cells.Item(50, 5).Value = "  -7.03%  "
$ws.Cells.Item(51, 4).Value = "'0.288"
$ws.Cells.Item(51, 5).Value = "  -7.04%  "
